$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) to reflect the reordered metric columns.
$ws.Range("B1").Value = "acc"
$ws.Range("C1").Value = "loss"
$ws.Range("D1").Value = "mean_absolute_error"
$ws.Range("E1").Value = "val_acc"
$ws.Range("F1").Value = "val_loss"
$ws.Range("G1").Value = "val_mean_absolute_error"

# Update data rows (2-41) with the new metric values per column.
$data = @(
    @(0.6485260795641378, 0.2339363681013081, 0.2339363681013081, 0.7328042302812848, 0.1964499006668727, 0.1964499006668727),
    @(0.7244897995000523, 0.1776597141003122, 0.1776597141003122, 0.7645502632887906, 0.142791933049916, 0.142791933049916),
    @(0.7823129222641726, 0.1254619428478266, 0.1254619428478266, 0.8544973519744066, 0.08517111991606062, 0.08517111991606062),
    @(0.8276643965249699, 0.08986162239698325, 0.08986162239698325, 0.8597883569500434, 0.05679801029581872, 0.05679801029581872),
    @(0.8605442121456, 0.07741538746928682, 0.07741538746928682, 0.8544973488207217, 0.04942771380461713, 0.04942771380461713),
    @(0.8458049856886571, 0.07219248730478103, 0.07219248730478103, 0.8730158686007141, 0.04670578676950995, 0.04670578676950995),
    @(0.8412698395127882, 0.06853477337530681, 0.06853477337530681, 0.8624338570725981, 0.04673912156392973, 0.04673912156392973),
    @(0.8628117870581664, 0.06883696155821115, 0.06883696155821115, 0.8624338586494406, 0.04869715664437208, 0.04869715664437208),
    @(0.8458049833233935, 0.06833576908558946, 0.06833576908558946, 0.8597883553732009, 0.04528263484241155, 0.04528263484241155),
    @(0.8684807226501084, 0.06722350878949339, 0.06722350878949339, 0.8597883569500434, 0.04688214109569946, 0.04688214109569946),
    @(0.8594104270545804, 0.06408613091524766, 0.06408613091524766, 0.8624338586494406, 0.05066828141925196, 0.05066828141925196),
    @(0.8673469308011926, 0.0645740676476977, 0.0645740676476977, 0.8597883569500434, 0.04672135800004951, 0.04672135800004951),
    @(0.8594104250272115, 0.06317227325336733, 0.06317227325336733, 0.8597883553732009, 0.04516204090779105, 0.04516204090779105),
    @(0.8752834426572804, 0.06278287619352341, 0.06278287619352341, 0.8756613718769538, 0.04462548879482759, 0.04462548879482759),
    @(0.8650793608894695, 0.06387534633684321, 0.06387534633684321, 0.8677248652019198, 0.04485266377765035, 0.04485266377765035),
    @(0.8571428507904347, 0.06122600532595127, 0.06122600532595127, 0.8650793635025227, 0.04338534040346978, 0.04338534040346978),
    @(0.8616780006155675, 0.06125166994475183, 0.06125166994475183, 0.8650793635025227, 0.04556353723324796, 0.04556353723324796),
    @(0.8639455735007652, 0.06037402545817856, 0.06037402545817856, 0.862433860226283, 0.04409303706356142, 0.04409303706356142),
    @(0.8730158670689243, 0.05889587356052161, 0.05889587356052161, 0.8571428552506462, 0.04559907312233927, 0.04559907312233927),
    @(0.8662131443585939, 0.05848183362087966, 0.05848183362087966, 0.862433860226283, 0.04671717172971478, 0.04671717172971478),
    @(0.8696145067274436, 0.05845955854742164, 0.05845955854742164, 0.8624338586494406, 0.04335810169183388, 0.04335810169183388),
    @(0.8786848015796029, 0.05910817637634115, 0.05910817637634115, 0.8571428552506462, 0.04385818069968274, 0.04385818069968274),
    @(0.8662131457101732, 0.05932440899751759, 0.05932440899751759, 0.8650793619256801, 0.04331020689594052, 0.04331020689594052),
    @(0.8809523741944847, 0.05623901362558341, 0.05623901362558341, 0.862433860226283, 0.04400462523181602, 0.04400462523181602),
    @(0.8809523741944847, 0.0570832078182508, 0.0570832078182508, 0.8597883569500434, 0.0433283896436767, 0.0433283896436767),
    @(0.8662131484133316, 0.05805420705757174, 0.05805420705757174, 0.862433860226283, 0.0440658430258433, 0.0440658430258433),
    @(0.8718820804911677, 0.05873115098023631, 0.05873115098023631, 0.8677248652019198, 0.0437937399579419, 0.0437937399579419),
    @(0.8707482935079371, 0.05624305710297863, 0.05624305710297863, 0.862433860226283, 0.04509334475117386, 0.04509334475117386),
    @(0.8764172263967207, 0.05858811671062121, 0.05858811671062121, 0.862433860226283, 0.04311702244732746, 0.04311702244732746),
    @(0.8730158680150298, 0.05630175254039483, 0.05630175254039483, 0.862433860226283, 0.04890885344021535, 0.04890885344021535),
    @(0.8798185870760963, 0.05527254869810849, 0.05527254869810849, 0.8677248652019198, 0.04301970048004358, 0.04301970048004358),
    @(0.8730158686908194, 0.05643147025721954, 0.05643147025721954, 0.8677248652019198, 0.04335480003997132, 0.04335480003997132),
    @(0.875283437926753, 0.05549544602835259, 0.05549544602835259, 0.8730158686007141, 0.04462635708312509, 0.04462635708312509),
    @(0.8820861606370836, 0.05622636406025649, 0.05622636406025649, 0.8650793619256801, 0.04771522864226311, 0.04771522864226311),
    @(0.8786848026608663, 0.05475998533015349, 0.05475998533015349, 0.870370366901317, 0.0431286519363759, 0.0431286519363759),
    @(0.8775510158127938, 0.05525013211430335, 0.05525013211430335, 0.8571428552506462, 0.04425570584597096, 0.04425570584597096),
    @(0.8741496548631024, 0.05579157561367871, 0.05579157561367871, 0.8650793619256801, 0.04256838499001725, 0.04256838499001725),
    @(0.8832199511344201, 0.05553818218697226, 0.05553818218697226, 0.862433860226283, 0.04280567773317218, 0.04280567773317218),
    @(0.8798185859948329, 0.05361374588408708, 0.05361374588408708, 0.862433860226283, 0.04301337968735468, 0.04301337968735468),
    @(0.8798185870760963, 0.05503128954714118, 0.05503128954714118, 0.8597883569500434, 0.04389099790542214, 0.04389099790542214)
)

$startRow = 2
$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($r, $j + 2).Value = $rowValues[$j]
    }
}
